$wb = $excel.ActiveWorkbook

# ---- Metadata sheet ----
$meta = $wb.Worksheets.Item("Metadata")

# Version bump
$meta.Range("B3").Value = "6.0.0"

# Date bump
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Row 10 was "Contact" / "No display for ContactDetail" -> becomes "Jurisdiction" / "United States of America"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row -> remove it entirely
$meta.Rows.Item(11).Delete()

# ---- Elements sheet ----
$elem = $wb.Worksheets.Item("Elements")

# Root Extension element row (row 2): Short & Definition now reflect this profile's title/description
$elem.Range("K2").Value = "Short Term Disability Benefit Rate"
$elem.Range("L2").Value = "Short term disability (STD) benefit rate, expressed as a percentage of base wages or benefit amount (for example, 50.00%, 60.00%, 70.00%)"
